# Auto-generated market-data refresh for the Masamune_Profits workbook.
# For each (sheet, row) pair below we overwrite the fetched-price columns
# (H/I/J/K/L = currentAveragePrice*, LevePrice*) and the derived profit
# columns (M/N = LeveProfit*) with the values from the latest scheduled
# market-data pull. When a row has no HQ/NQ listings at all, the profit
# cell for that side is left blank (cleared) rather than holding a stale
# number, matching how the rest of the sheet represents "no data".

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 732.875
$ws.Range("I28").Value = 644.1667
$ws.Range("K28").Value = 644.1667
$ws.Range("M28").Value = -159.1667

# Row 95
$ws.Range("H95").Value = 33330.668
$ws.Range("J95").Value = 33330.668
$ws.Range("L95").Value = 33330.668
$ws.Range("N95").Value = -38822.668

# Row 105
$ws.Range("H105").Value = 48831.5
$ws.Range("J105").Value = 48831.5
$ws.Range("L105").Value = 48831.5
$ws.Range("N105").Value = -55819.5

$ws = $wb.Worksheets.Item("ARM")
# Row 7
$ws.Range("H7").Value = 50000
$ws.Range("J7").Value = 50000
$ws.Range("L7").Value = 50000
$ws.Range("N7").Value = -50228

# Row 95
$ws.Range("H95").Value = 35396.8
$ws.Range("J95").Value = 35396.8
$ws.Range("L95").Value = 35396.8
$ws.Range("N95").Value = -40888.8

# Row 101
$ws.Range("H101").Value = 44496.8
$ws.Range("J101").Value = 44496.8
$ws.Range("L101").Value = 44496.8
$ws.Range("N101").Value = -50986.8

# Row 103
$ws.Range("H103").Value = 37996
$ws.Range("J103").Value = 37996
$ws.Range("L103").Value = 37996
$ws.Range("N103").Value = -40340

# Row 104
$ws.Range("H104").Value = 40733.668
$ws.Range("J104").Value = 40733.668
$ws.Range("L104").Value = 40733.668
$ws.Range("N104").Value = -47721.668

# Row 105
$ws.Range("H105").Value = 38389.2
$ws.Range("J105").Value = 38389.2
$ws.Range("L105").Value = 38389.2
$ws.Range("N105").Value = -45377.2

# Row 106
$ws.Range("H106").Value = 46338.5
$ws.Range("J106").Value = 46338.5
$ws.Range("L106").Value = 46338.5
$ws.Range("N106").Value = -48862.5

# Row 121
$ws.Range("H121").Value = 33773
$ws.Range("J121").Value = 33773
$ws.Range("L121").Value = 33773
$ws.Range("N121").Value = -37267

# Row 123
$ws.Range("H123").Value = 40878
$ws.Range("J123").Value = 40878
$ws.Range("L123").Value = 40878
$ws.Range("N123").Value = -50678

$ws = $wb.Worksheets.Item("BSM")
# Row 95
$ws.Range("H95").Value = 40997.332
$ws.Range("J95").Value = 40997.332
$ws.Range("L95").Value = 40997.332
$ws.Range("N95").Value = -46489.332

# Row 100
$ws.Range("H100").Value = 9513.166999999999
$ws.Range("J100").Value = 9513.166999999999
$ws.Range("L100").Value = 9513.166999999999
$ws.Range("N100").Value = -11677.167

# Row 103
$ws.Range("H103").Value = 40185.145
$ws.Range("J103").Value = 40185.145
$ws.Range("L103").Value = 40185.145
$ws.Range("N103").Value = -42529.145

# Row 122
$ws.Range("H122").Value = 36079
$ws.Range("J122").Value = 36079
$ws.Range("L122").Value = 36079
$ws.Range("N122").Value = -45879

$ws = $wb.Worksheets.Item("CRP")
# Row 28
$ws.Range("H28").Value = 32555.715
$ws.Range("J28").Value = 32555.715
$ws.Range("L28").Value = 32555.715
$ws.Range("N28").Value = -33045.715

# Row 43
$ws.Range("H43").Value = 38885.668
$ws.Range("J43").Value = 38885.668
$ws.Range("L43").Value = 38885.668
$ws.Range("N43").Value = -39253.668

# Row 101
$ws.Range("H101").Value = 38885.668
$ws.Range("J101").Value = 38885.668
$ws.Range("L101").Value = 38885.668
$ws.Range("N101").Value = -45375.668

# Row 106
$ws.Range("H106").Value = 32635.334
$ws.Range("J106").Value = 32635.334
$ws.Range("L106").Value = 32635.334
$ws.Range("N106").Value = -35159.334

# Row 111
$ws.Range("H111").Value = 26684
$ws.Range("J111").Value = 26684
$ws.Range("L111").Value = 26684
$ws.Range("N111").Value = -34864

# Row 119
$ws.Range("H119").Value = 35487.332
$ws.Range("J119").Value = 35487.332
$ws.Range("L119").Value = 35487.332
$ws.Range("N119").Value = -45163.332

# Row 132
$ws.Range("H132").Value = 54078.297
$ws.Range("I132").Value = 1700.25
$ws.Range("K132").Value = 5100.75
$ws.Range("M132").Value = -2570.75

$ws = $wb.Worksheets.Item("CUL")
# Row 126
$ws.Range("H126").Value = 26597.924
$ws.Range("I126").Value = 300000
$ws.Range("J126").Value = 3814.4167
$ws.Range("K126").Value = 900000
$ws.Range("L126").Value = 11443.2501
$ws.Range("M126").Value = -895060
$ws.Range("N126").Value = -21323.2501

$ws = $wb.Worksheets.Item("GSM")
# Row 6
$ws.Range("H6").Value = 20898
$ws.Range("J6").Value = 19872.5
$ws.Range("L6").Value = 19872.5
$ws.Range("N6").Value = -20098.5

# Row 9
$ws.Range("H9").Value = 2450
$ws.Range("I9").Value = 2450
$ws.Range("K9").Value = 2450
$ws.Range("M9").Value = -2280

# Row 12
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()

# Row 16
$ws.Range("H16").Value = 20898
$ws.Range("J16").Value = 19872.5
$ws.Range("L16").Value = 19872.5
$ws.Range("N16").Value = -20372.5

# Row 98
$ws.Range("H98").Value = 33776.285
$ws.Range("J98").Value = 33776.285
$ws.Range("L98").Value = 33776.285
$ws.Range("N98").Value = -39766.285

# Row 101
$ws.Range("H101").Value = 40653
$ws.Range("J101").Value = 40653
$ws.Range("L101").Value = 40653
$ws.Range("N101").Value = -47143

# Row 104
$ws.Range("H104").Value = 33797.4
$ws.Range("J104").Value = 33797.4
$ws.Range("L104").Value = 33797.4
$ws.Range("N104").Value = -40785.4

# Row 105
$ws.Range("H105").Value = 38521.145
$ws.Range("J105").Value = 38521.145
$ws.Range("L105").Value = 38521.145
$ws.Range("N105").Value = -45509.145

# Row 110
$ws.Range("H110").Value = 34628.75
$ws.Range("J110").Value = 34628.75
$ws.Range("L110").Value = 34628.75
$ws.Range("N110").Value = -42808.75

$ws = $wb.Worksheets.Item("LTW")
# Row 10
$ws.Range("H10").Value = 2000
$ws.Range("J10").Value = 2000
$ws.Range("L10").Value = 2000
$ws.Range("N10").Value = -2280

# Row 103
$ws.Range("H103").Value = 41187.332
$ws.Range("J103").Value = 41187.332
$ws.Range("L103").Value = 41187.332
$ws.Range("N103").Value = -43531.332

# Row 105
$ws.Range("H105").Value = 33579.625
$ws.Range("J105").Value = 33579.625
$ws.Range("L105").Value = 33579.625
$ws.Range("N105").Value = -40567.625

# Row 106
$ws.Range("H106").Value = 31994
$ws.Range("J106").Value = 31994
$ws.Range("L106").Value = 31994
$ws.Range("N106").Value = -34518

# Row 110
$ws.Range("H110").Value = 35817
$ws.Range("J110").Value = 35817
$ws.Range("L110").Value = 35817
$ws.Range("N110").Value = -43997

# Row 114
$ws.Range("H114").Value = 23567.666
$ws.Range("J114").Value = 23567.666
$ws.Range("L114").Value = 23567.666
$ws.Range("N114").Value = -32245.666

# Row 120
$ws.Range("H120").Value = 37437.25
$ws.Range("J120").Value = 37437.25
$ws.Range("L120").Value = 37437.25
$ws.Range("N120").Value = -47113.25

# Row 121
$ws.Range("H121").Value = 21259.334
$ws.Range("J121").Value = 21259.334
$ws.Range("L121").Value = 21259.334
$ws.Range("N121").Value = -24753.334

$ws = $wb.Worksheets.Item("WVR")
# Row 95
$ws.Range("H95").Value = 39988
$ws.Range("J95").Value = 39988
$ws.Range("L95").Value = 39988
$ws.Range("N95").Value = -45480

# Row 97
$ws.Range("H97").Value = 35693
$ws.Range("J97").Value = 35693
$ws.Range("L97").Value = 35693
$ws.Range("N97").Value = -37675

# Row 98
$ws.Range("H98").Value = 28346.666
$ws.Range("J98").Value = 28346.666
$ws.Range("L98").Value = 28346.666
$ws.Range("N98").Value = -34336.666

# Row 103
$ws.Range("H103").Value = 35722
$ws.Range("J103").Value = 35722
$ws.Range("L103").Value = 35722
$ws.Range("N103").Value = -38066

# Row 104
$ws.Range("H104").Value = 39344.4
$ws.Range("J104").Value = 39344.4
$ws.Range("L104").Value = 39344.4
$ws.Range("N104").Value = -46332.4

# Row 106
$ws.Range("H106").Value = 32479.777
$ws.Range("J106").Value = 32479.777
$ws.Range("L106").Value = 32479.777
$ws.Range("N106").Value = -35003.777

# Row 110
$ws.Range("H110").Value = 26619.75
$ws.Range("J110").Value = 26619.75
$ws.Range("L110").Value = 26619.75
$ws.Range("N110").Value = -34799.75

# Row 112
$ws.Range("H112").Value = 26375
$ws.Range("J112").Value = 26375
$ws.Range("L112").Value = 26375
$ws.Range("N112").Value = -29329

# Row 116
$ws.Range("H116").Value = 27061.5
$ws.Range("J116").Value = 27061.5
$ws.Range("L116").Value = 27061.5
$ws.Range("N116").Value = -36239.5

# Row 118
$ws.Range("H118").Value = 25278.4
$ws.Range("J118").Value = 27848
$ws.Range("L118").Value = 27848
$ws.Range("N118").Value = -31162
